$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test case row (row 20) - "Customercare017"
# Shared-string insertion order matters for index assignment (50,51,52),
# so write in the same order as the target shared strings table:
# 50 = Description (C), 51 = Jira id (B), 52 = TCID (A)
$ws.Range("C20").Value = "Verify that Countries list should be updated to match with SFDC list in customer care page as per document OPWLRA-630.xlsx."
$ws.Range("B20").Value = "OPQA-5350"
$ws.Range("A20").Value = "Customercare017"
$ws.Range("D20").Value = "Y"

# Update selection to match new active cell
$ws.Range("B20").Select()

# Formatting tweaks on row 20 to match target style indices
$ws.Range("B20").WrapText = $true
$ws.Range("D19").Copy()
$ws.Range("D20").PasteSpecial(-4122)
